$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 918.8
$ws.Range("I61").Value = 698
$ws.Range("J61").Value = 1250
$ws.Range("K61").Value = 2094
$ws.Range("L61").Value = 3750
$ws.Range("M61").Value = -1922
$ws.Range("N61").Value = -4094

$ws.Range("H76").Value = 7582
$ws.Range("I76").Value = 7065
$ws.Range("K76").Value = 7065
$ws.Range("M76").Value = -6750

$ws.Range("H79").Value = 7582
$ws.Range("I79").Value = 7065
$ws.Range("K79").Value = 7065
$ws.Range("M79").Value = -5973

$ws.Range("H100").Value = 1836.2
$ws.Range("I100").Value = 1111.625
$ws.Range("J100").Value = 4734.5
$ws.Range("K100").Value = 1111.625
$ws.Range("L100").Value = 4734.5
$ws.Range("M100").Value = -570.625
$ws.Range("N100").Value = -5816.5

$ws.Range("H135").Value = 2000459
$ws.Range("I135").Value = 2500373.8
$ws.Range("K135").Value = 22503364.2
$ws.Range("M135").Value = -22500829.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 24970
$ws.Range("J42").Value = 24970
$ws.Range("L42").Value = 24970
$ws.Range("N42").Value = -25942

$ws.Range("H97").Value = 3206372.5
$ws.Range("I97").Value = 1097.4584
$ws.Range("J97").Value = 41669670
$ws.Range("K97").Value = 1097.4584
$ws.Range("L97").Value = 41669670
$ws.Range("M97").Value = -601.4584
$ws.Range("N97").Value = -41670662

$ws.Range("H123").Value = 58849
$ws.Range("J123").Value = 58849
$ws.Range("L123").Value = 58849
$ws.Range("N123").Value = -68649

$ws.Range("H132").Value = 7350.2925
$ws.Range("I132").Value = 5429.5
$ws.Range("K132").Value = 16288.5
$ws.Range("M132").Value = -13758.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 102790.5
$ws.Range("I86").Value = 251051.75
$ws.Range("J86").Value = 3949.6667
$ws.Range("K86").Value = 251051.75
$ws.Range("L86").Value = 3949.6667
$ws.Range("M86").Value = -249928.75
$ws.Range("N86").Value = -6195.6667

$ws.Range("H89").Value = 102790.5
$ws.Range("I89").Value = 251051.75
$ws.Range("J89").Value = 3949.6667
$ws.Range("K89").Value = 1255258.75
$ws.Range("L89").Value = 19748.3335
$ws.Range("M89").Value = -1249642.75
$ws.Range("N89").Value = -30980.3335

$ws.Range("H94").Value = 1676.2273
$ws.Range("I94").Value = 1093.6
$ws.Range("J94").Value = 7502.5
$ws.Range("K94").Value = 1093.6
$ws.Range("L94").Value = 7502.5
$ws.Range("M94").Value = -642.5999999999999
$ws.Range("N94").Value = -8404.5

$ws.Range("H99").Value = 3014.2856
$ws.Range("J99").Value = 3096.2
$ws.Range("L99").Value = 3096.2
$ws.Range("N99").Value = -6092.2

$ws.Range("H134").Value = 5440217.5
$ws.Range("I134").Value = 11906388
$ws.Range("J134").Value = 8634.16
$ws.Range("K134").Value = 35719164
$ws.Range("L134").Value = 25902.48
$ws.Range("M134").Value = -35716629
$ws.Range("N134").Value = -30972.48

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7854.3237
$ws.Range("I31").Value = 2018.0625
$ws.Range("K31").Value = 2018.0625
$ws.Range("M31").Value = -1723.0625

$ws.Range("H34").Value = 7854.3237
$ws.Range("I34").Value = 2018.0625
$ws.Range("K34").Value = 2018.0625
$ws.Range("M34").Value = -1816.0625

$ws.Range("H96").Value = 42567.2
$ws.Range("J96").Value = 42567.2
$ws.Range("L96").Value = 42567.2
$ws.Range("N96").Value = -48059.2

$ws.Range("H107").Value = 2385.8635
$ws.Range("J107").Value = 2628.4211
$ws.Range("L107").Value = 2628.4211
$ws.Range("N107").Value = -6468.4211

$ws.Range("H132").Value = 6508.76
$ws.Range("I132").Value = 3594.1538
$ws.Range("J132").Value = 9666.25
$ws.Range("K132").Value = 10782.4614
$ws.Range("L132").Value = 28998.75
$ws.Range("M132").Value = -8252.4614
$ws.Range("N132").Value = -34058.75

$ws.Range("H134").Value = 6596.971
$ws.Range("J134").Value = 7944.2915
$ws.Range("L134").Value = 23832.8745
$ws.Range("N134").Value = -28902.8745

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 4999.3335
$ws.Range("I87").Value = 3999.5
$ws.Range("K87").Value = 11998.5
$ws.Range("M87").Value = -10750.5

$ws.Range("H90").Value = 4999.3335
$ws.Range("I90").Value = 3999.5
$ws.Range("K90").Value = 35995.5
$ws.Range("M90").Value = -29755.5

$ws.Range("H97").Value = 572
$ws.Range("I97").Value = 572
$ws.Range("K97").Value = 1716
$ws.Range("M97").Value = -1220

$ws.Range("H115").Value = 1195.8572
$ws.Range("I115").Value = 1195.8572
$ws.Range("K115").Value = 3587.5716
$ws.Range("M115").Value = -2412.5716

$ws.Range("H117").Value = 2109.375
$ws.Range("I117").Value = 1995
$ws.Range("J117").Value = 2147.5
$ws.Range("K117").Value = 5985
$ws.Range("L117").Value = 6442.5
$ws.Range("M117").Value = -2543
$ws.Range("N117").Value = -13326.5

$ws.Range("H128").Value = 109975
$ws.Range("I128").Value = 109975
$ws.Range("K128").Value = 329925
$ws.Range("M128").Value = -324945

$ws.Range("H130").Value = 3583.3333
$ws.Range("I130").Value = 2875
$ws.Range("K130").Value = 8625
$ws.Range("M130").Value = -3605

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 83749.75
$ws.Range("J52").Value = 90000
$ws.Range("L52").Value = 90000
$ws.Range("N52").Value = -90518

$ws.Range("H58").Value = 74719.60000000001
$ws.Range("J58").Value = 74719.60000000001
$ws.Range("L58").Value = 74719.60000000001
$ws.Range("N58").Value = -75273.60000000001

$ws.Range("H97").Value = 681.35
$ws.Range("I97").Value = 589.4643
$ws.Range("J97").Value = 895.75
$ws.Range("K97").Value = 589.4643
$ws.Range("L97").Value = 895.75
$ws.Range("M97").Value = -93.46429999999998
$ws.Range("N97").Value = -1887.75

$ws.Range("H102").Value = 3108.9666
$ws.Range("I102").Value = 2968.52
$ws.Range("K102").Value = 2968.52
$ws.Range("M102").Value = -1346.52

$ws.Range("H107").Value = 445054.22
$ws.Range("I107").Value = 728102.9399999999
$ws.Range("K107").Value = 728102.9399999999
$ws.Range("M107").Value = -726182.9399999999

$ws.Range("H122").Value = 4530023
$ws.Range("I122").Value = 6037531
$ws.Range("K122").Value = 18112593
$ws.Range("M122").Value = -18110143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4614.9565
$ws.Range("I7").Value = 3167.818
$ws.Range("K7").Value = 3167.818
$ws.Range("M7").Value = -3055.818

$ws.Range("H33").Value = 16669999
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H100").Value = 3559.182
$ws.Range("I100").Value = 1956.1333
$ws.Range("J100").Value = 6994.2856
$ws.Range("K100").Value = 1956.1333
$ws.Range("L100").Value = 6994.2856
$ws.Range("M100").Value = -1415.1333
$ws.Range("N100").Value = -8076.2856

$ws.Range("H126").Value = 4614.9565
$ws.Range("I126").Value = 3167.818
$ws.Range("K126").Value = 9503.454000000002
$ws.Range("M126").Value = -7033.454000000002

$ws.Range("H136").Value = 12185.511
$ws.Range("I136").Value = 2829.923
$ws.Range("K136").Value = 8489.769
$ws.Range("M136").Value = -5939.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 28405
$ws.Range("J41").Value = 28685
$ws.Range("L41").Value = 28685
$ws.Range("N41").Value = -29465

$ws.Range("H122").Value = 104632.64
$ws.Range("I122").Value = 134866.7
$ws.Range("J122").Value = 3852.4443
$ws.Range("K122").Value = 404600.1
$ws.Range("L122").Value = 11557.3329
$ws.Range("M122").Value = -402150.1
$ws.Range("N122").Value = -16457.3329

$ws.Range("H123").Value = 50567
$ws.Range("J123").Value = 50567
$ws.Range("L123").Value = 50567
$ws.Range("N123").Value = -60367

$ws.Range("H130").Value = 74759.5
$ws.Range("J130").Value = 74759.5
$ws.Range("L130").Value = 74759.5
$ws.Range("N130").Value = -84799.5
